$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Range("D4").Value = "l10n_it_ricevute_bancarie"
$ws.Range("A19").Select()
